# Apply TradeXCB strategy sheet updates (row 2 = Buy leg, row 3 = Sell leg).
# Leading apostrophes force Excel to store the values as text (matching the
# original inlineStr cell type) instead of auto-coercing to numbers/dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Buy leg)
$ws.Range("I2").Value = "'2022-04-21"
$ws.Range("J2").Value = "'NIFTY2242117500CE"
$ws.Range("K2").Value = "'130"
$ws.Range("L2").Value = "'YES"
$ws.Range("X2").Value = "'NO"

# Row 3 (Sell leg)
$ws.Range("I3").Value = "'2022-04-21"
$ws.Range("J3").Value = "'NIFTY2242117500PE"
$ws.Range("M3").Value = "'150"
$ws.Range("N3").Value = "'YES"
$ws.Range("X3").Value = "'NO"

Write-Host "Applied TradeXCB strategy updates"
